$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.171.57"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.703.34"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("D5").Value = "'223.58"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'0.5235"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "'0.06593"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.2617"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").Value = "'20.56"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "'0.07731"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.423"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "1.939.45"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.702.46"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'0.5731"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "0.0₅8102"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "'66.95"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "27.274.04"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "217.13"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'1.009"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'4.595"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").Value = "10.31"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "'5.987"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").Value = "144.88"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "1.728"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "'0.1194"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "'7.156"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").Value = "'16.00"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "'0.05273"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "'1.290"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "3.423"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "'3.307"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").Value = "1.623"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'2.815"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.405"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9419"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "'0.5818"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "1.180.67"
$ws.Range("E39").Value = "  +13.40%  "
$ws.Range("D40").Value = "'0.01633"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "'1.010"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'5.749"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "'0.8361"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").Value = "100.77"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "1.851.20"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("D47").Value = "'57.00"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "'0.4557"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "'8.053"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'0.05226"
$ws.Range("E51").Value = "  -0.28%  "
